# "Generate Report for Handback" - update localization-status workbook
# with the results of a handback run for zh-cn and de-de.

$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$statusHandedBack = "Handed back: in sync with en-US"

# --- Status updates: "Ready for handoff" -> "Handed back: in sync with en-US" ---
$wsOverview.Range("E2").Value = $statusHandedBack
$wsOverview.Range("F2").Value = $statusHandedBack
$wsOverview.Range("E3").Value = $statusHandedBack
$wsOverview.Range("F3").Value = $statusHandedBack

$wsZhCn.Range("C2").Value = $statusHandedBack
$wsZhCn.Range("C3").Value = $statusHandedBack

$wsDeDe.Range("C2").Value = $statusHandedBack
$wsDeDe.Range("C3").Value = $statusHandedBack

# --- zh-cn handback details ---
$wsZhCn.Range("I2").Value = "6dc7918b-f1eb-4e82-8597-1242ab6771e1.md"
$wsZhCn.Range("I2").Style = "HyperLink"
$wsZhCn.Range("J2").Value = "6dc7918b-f1eb-4e82-8597-1242ab6771e1.fb1cedfe47768b6baef2737f81e3c7cf9f28dcf9.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-31 08:59:24"

$wsZhCn.Range("I3").Value = "ae31b1f1-136c-4935-a87b-de54edff9fbe.md"
$wsZhCn.Range("I3").Style = "HyperLink"
$wsZhCn.Range("J3").Value = "ae31b1f1-136c-4935-a87b-de54edff9fbe.08fbeca7bfcb5d5c5f6e7810fa72756d776553d9.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-31 08:59:24"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b96bdb51f10f73bd4dd833bf26240ffb352c0239/e2e/6dc7918b-f1eb-4e82-8597-1242ab6771e1.md", "", "", "6dc7918b-f1eb-4e82-8597-1242ab6771e1.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b96bdb51f10f73bd4dd833bf26240ffb352c0239/e2e/ae31b1f1-136c-4935-a87b-de54edff9fbe.md", "", "", "ae31b1f1-136c-4935-a87b-de54edff9fbe.md")

# --- de-de handback details ---
$wsDeDe.Range("I2").Value = "6dc7918b-f1eb-4e82-8597-1242ab6771e1.md"
$wsDeDe.Range("I2").Style = "HyperLink"
$wsDeDe.Range("J2").Value = "6dc7918b-f1eb-4e82-8597-1242ab6771e1.fb1cedfe47768b6baef2737f81e3c7cf9f28dcf9.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-31 08:59:44"

$wsDeDe.Range("I3").Value = "ae31b1f1-136c-4935-a87b-de54edff9fbe.md"
$wsDeDe.Range("I3").Style = "HyperLink"
$wsDeDe.Range("J3").Value = "ae31b1f1-136c-4935-a87b-de54edff9fbe.08fbeca7bfcb5d5c5f6e7810fa72756d776553d9.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-31 08:59:44"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b96bdb51f10f73bd4dd833bf26240ffb352c0239/e2e/6dc7918b-f1eb-4e82-8597-1242ab6771e1.md", "", "", "6dc7918b-f1eb-4e82-8597-1242ab6771e1.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b96bdb51f10f73bd4dd833bf26240ffb352c0239/e2e/ae31b1f1-136c-4935-a87b-de54edff9fbe.md", "", "", "ae31b1f1-136c-4935-a87b-de54edff9fbe.md")

# --- Column widths: widen columns that now hold longer handback text ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.17
$wsOverview.Columns.Item(6).ColumnWidth = 29.17

$wsZhCn.Columns.Item(3).ColumnWidth = 29.17
$wsZhCn.Columns.Item(9).ColumnWidth = 39.17
$wsZhCn.Columns.Item(10).ColumnWidth = 39.17

$wsDeDe.Columns.Item(3).ColumnWidth = 29.17
$wsDeDe.Columns.Item(9).ColumnWidth = 39.17
$wsDeDe.Columns.Item(10).ColumnWidth = 39.17
